# Fruta / hortaliza, semanal
# Insert a new weekly record at row 32 (pushing the existing rows 32-55
# down to 33-56) on the Chirimoya / Macroferia Regional de Talca sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 32:55 down to 33:56, leaving a blank row 32 to fill in.
$ws.Rows("32:32").Insert()

$ws.Range("A32").Value = 5
$ws.Range("B32").Value = "Macroferia Regional de Talca"
$ws.Range("C32").Value = "Maule"
$ws.Range("D32").Value = 44512
$ws.Range("E32").Value = 7
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100107
$ws.Range("H32").Value = "Otros"
$ws.Range("I32").Value = 100107002
$ws.Range("J32").Value = "Chirimoya"
$ws.Range("K32").Value = "Cultivar IV Región"
$ws.Range("L32").Value = "Especial"
$ws.Range("M32").Value = 150
$ws.Range("N32").Value = 26000
$ws.Range("O32").Value = 26000
$ws.Range("P32").Value = 26000
$ws.Range("Q32").Value = "$/bandeja 10 kilos"
$ws.Range("R32").Value = "Provincia de Limarí"
$ws.Range("S32").Value = 2600
$ws.Range("T32").Value = 10
